# Applies the "Added test scripts in to IAM" change: appends rows 20-27
# to the IAM test-data sheet (sheet1 / "IAM"), extending the used range
# from A1:L19 to A1:L27 and updating shared strings accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("A20").Value = "OPQA-745"
$ws.Range("B20").Value = "Verify that to evict not-logged-in user and test whether that user is able to login or not"
$ws.Range("C20").Value = "1PAUTH"
$ws.Range("D20").Value = "/admin/access"
$ws.Range("D20").Style = "Hyperlink"
$ws.Range("E20").Value = "PUT"
$ws.Range("F20").Value = "Content-Type=application/json"
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").Value = "{`"truid`":`"(SYS_USER2)`" ,`"userStatus`":`"Deactivate`", `"comments`":`"Deactivate User`"}"
$ws.Range("J20").Value = "status=200"
$ws.Rows.Item(20).RowHeight = 45

# Row 21
$ws.Range("A21").Value = "OPQA-543"
$ws.Range("B21").Value = "Verify that evicted user  not able to login with valid credentials "
$ws.Range("C21").Value = "1PAUTH"
$ws.Range("D21").Value = "/authorize"
$ws.Range("D21").Style = "Hyperlink"
$ws.Range("E21").Value = "POST"
$ws.Range("F21").Value = "Content-Type=application/json"
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value = "{`"loginid`":`"(OPQA-542_email)`",`"password`":`"1Platform!`"}"
$ws.Range("I21").Value = "OPQA-542"
$ws.Range("J21").Value = "status=423||reason=User is evicted||errorcode=423"
$ws.Rows.Item(21).RowHeight = 30

# Row 22
$ws.Range("A22").Value = "OPQA-746"
$ws.Range("B22").Value = "Verify that user is able to activate evicted user by passing truid"
$ws.Range("C22").Value = "1PAUTH"
$ws.Range("D22").Value = "/admin/access"
$ws.Range("D22").Style = "Hyperlink"
$ws.Range("E22").Value = "PUT"
$ws.Range("F22").Value = "Content-Type=application/json"
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").Value = "{`"truid`":`"(SYS_USER2)`" ,`"userStatus`":`"Activate`", `"comments`":`"Activate evicted User`"}"
$ws.Range("J22").Value = "status=200"
$ws.Rows.Item(22).RowHeight = 30

# Row 23
$ws.Range("A23").Value = "OPQA-2706_1"
$ws.Range("B23").Value = "Verify that user able to login with valid credentials"
$ws.Range("C23").Value = "1PAUTH"
$ws.Range("D23").Value = "/authorize"
$ws.Range("D23").Style = "Hyperlink"
$ws.Range("E23").Value = "POST"
$ws.Range("F23").Value = "Content-Type=application/json"
$ws.Range("G23").Style = "Normal"
$ws.Range("H23").Value = "{`"loginid`":`"(OPQA-542_email)`",`"password`":`"1Platform!`"}"
$ws.Range("I23").Value = "OPQA-542"
$ws.Range("J23").Value = "status=200||userid=(SYS_USER2)"
$ws.Range("K23").Value = "token"
$ws.Rows.Item(23).RowHeight = 30

# Row 24
$ws.Range("A24").Value = "OPQA-544"
$ws.Range("B24").Value = "Verify that system is ability to evict an already logged in user. "
$ws.Range("C24").Value = "1PAUTH"
$ws.Range("D24").Value = "/admin/access"
$ws.Range("D24").Style = "Hyperlink"
$ws.Range("E24").Value = "PUT"
$ws.Range("F24").Value = "Content-Type=application/json"
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").Value = "{`"truid`":`"(SYS_USER2)`" ,`"userStatus`":`"Deactivate`", `"comments`":`"Deactivate User`"}"
$ws.Range("J24").Value = "status=200"
$ws.Rows.Item(24).RowHeight = 30

# Row 25
$ws.Range("A25").Value = "OPQA-544_1"
$ws.Range("B25").Value = "Verify that to validate token which was evicted user who already login to neon "
$ws.Range("C25").Value = "1PAUTH"
$ws.Range("D25").Value = "/validate/(OPQA-2706_1_token)"
$ws.Range("D25").Style = "Hyperlink"
$ws.Range("E25").Value = "GET"
$ws.Range("G25").Style = "Normal"
$ws.Range("H25").Style = "Normal"
$ws.Range("I25").Value = "OPQA-2706_1"
$ws.Range("J25").Value = "status=423"
$ws.Rows.Item(25).RowHeight = 45

# Row 26
$ws.Range("A26").Value = "OPQA-544_2"
$ws.Range("B26").Value = "Verify that user is able to activate evicted user by passing truid"
$ws.Range("C26").Value = "1PAUTH"
$ws.Range("D26").Value = "/admin/access"
$ws.Range("D26").Style = "Hyperlink"
$ws.Range("E26").Value = "PUT"
$ws.Range("F26").Value = "Content-Type=application/json"
$ws.Range("G26").Style = "Normal"
$ws.Range("H26").Value = "{`"truid`":`"(SYS_USER2)`" ,`"userStatus`":`"Activate`", `"comments`":`"Activate evicted User`"}"
$ws.Range("J26").Value = "status=200"
$ws.Rows.Item(26).RowHeight = 30

# Row 27
$ws.Range("A27").Value = "OPQA-547"
$ws.Range("B27").Value = "Verify whether reverted user able to log in or not. And test reverted user should able to login in to Neon"
$ws.Range("C27").Value = "1PAUTH"
$ws.Range("D27").Value = "/authorize"
$ws.Range("D27").Style = "Hyperlink"
$ws.Range("E27").Value = "POST"
$ws.Range("F27").Value = "Content-Type=application/json"
$ws.Range("G27").Style = "Normal"
$ws.Range("H27").Value = "{`"loginid`":`"(OPQA-542_email)`",`"password`":`"1Platform!`"}"
$ws.Range("I27").Value = "OPQA-542"
$ws.Range("J27").Value = "status=200||userid=(SYS_USER2)"
$ws.Rows.Item(27).RowHeight = 60

# Update the sheet selection to mirror the new extent of the data
$ws.Range("L2:L27").Select()
